# Insert 3 new price rows (date 44706) for "Packham's Triumph" at the top of
# this block of weekly Pera (pear) price records, pushing the existing
# records down by 3 rows (dimension grows from A1:T1255 to A1:T1258).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 1177:1255 down to 1180:1258, keeping their formatting.
$ws.Rows("1177:1179").Insert()

# Row 1177: Packham's Triumph / Especial
$ws.Range("A1177").Value = 3
$ws.Range("B1177").Value = "Femacal de La Calera"
$ws.Range("C1177").Value = "Coquimbo"
$ws.Range("D1177").Value = 44706
$ws.Range("E1177").Value = 5
$ws.Range("F1177").Value = "Fruta"
$ws.Range("G1177").Value = 100104
$ws.Range("H1177").Value = "Frutos de pepita"
$ws.Range("I1177").Value = 100104005
$ws.Range("J1177").Value = "Pera"
$ws.Range("K1177").Value = "Packham's Triumph"
$ws.Range("L1177").Value = "Especial"
$ws.Range("M1177").Value = 75
$ws.Range("N1177").Value = 10000
$ws.Range("O1177").Value = 10000
$ws.Range("P1177").Value = 10000
$ws.Range("Q1177").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R1177").Value = "Región de O'Higgins"
$ws.Range("S1177").Value = 556
$ws.Range("T1177").Value = 18

# Row 1178: Packham's Triumph / Primera
$ws.Range("A1178").Value = 3
$ws.Range("B1178").Value = "Femacal de La Calera"
$ws.Range("C1178").Value = "Coquimbo"
$ws.Range("D1178").Value = 44706
$ws.Range("E1178").Value = 5
$ws.Range("F1178").Value = "Fruta"
$ws.Range("G1178").Value = 100104
$ws.Range("H1178").Value = "Frutos de pepita"
$ws.Range("I1178").Value = 100104005
$ws.Range("J1178").Value = "Pera"
$ws.Range("K1178").Value = "Packham's Triumph"
$ws.Range("L1178").Value = "Primera"
$ws.Range("M1178").Value = 85
$ws.Range("N1178").Value = 9000
$ws.Range("O1178").Value = 9000
$ws.Range("P1178").Value = 9000
$ws.Range("Q1178").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R1178").Value = "Región de O'Higgins"
$ws.Range("S1178").Value = 500
$ws.Range("T1178").Value = 18

# Row 1179: Packham's Triumph / Segunda
$ws.Range("A1179").Value = 3
$ws.Range("B1179").Value = "Femacal de La Calera"
$ws.Range("C1179").Value = "Coquimbo"
$ws.Range("D1179").Value = 44706
$ws.Range("E1179").Value = 5
$ws.Range("F1179").Value = "Fruta"
$ws.Range("G1179").Value = 100104
$ws.Range("H1179").Value = "Frutos de pepita"
$ws.Range("I1179").Value = 100104005
$ws.Range("J1179").Value = "Pera"
$ws.Range("K1179").Value = "Packham's Triumph"
$ws.Range("L1179").Value = "Segunda"
$ws.Range("M1179").Value = 80
$ws.Range("N1179").Value = 8000
$ws.Range("O1179").Value = 8000
$ws.Range("P1179").Value = 8000
$ws.Range("Q1179").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R1179").Value = "Región de O'Higgins"
$ws.Range("S1179").Value = 444
$ws.Range("T1179").Value = 18

Write-Host "Inserted 3 new rows at 1177:1179; dimension now" $ws.UsedRange.Address()
